# Update Top50_DataComp.xlsx data: bump length counters (M2_Len) by 1 and
# roll forward the M2 date range (M2_1stDate / M2_LastDate) for the rows
# whose underlying data series were extended by one additional month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only the M2_Len (C) and M2_LastDate (F) change
$lenAndLastDateRows = @(
    @{ Row = 2;  C = 331; F = 45108 },
    @{ Row = 33; C = 451; F = 45108 },
    @{ Row = 34; C = 199; F = 45108 },
    @{ Row = 36; C = 499; F = 45108 },
    @{ Row = 40; C = 283; F = 45108 },
    @{ Row = 42; C = 235; F = 45078 },
    @{ Row = 44; C = 402; F = 45078 },
    @{ Row = 49; C = 294; F = 45078 },
    @{ Row = 50; C = 356; F = 45108 },
    @{ Row = 52; C = 343; F = 45108 }
)

foreach ($item in $lenAndLastDateRows) {
    $ws.Cells.Item($item.Row, 3).Value = $item.C   # Column C = M2_Len
    $ws.Cells.Item($item.Row, 6).Value = $item.F   # Column F = M2_LastDate
}

# Rows where the M2_1stDate (E) and M2_LastDate (F) change
$firstAndLastDateRows = @(
    @{ Row = 3;  E = 29921; F = 45108 },
    @{ Row = 5;  E = 29921; F = 45108 },
    @{ Row = 7;  E = 29891; F = 45078 },
    @{ Row = 11; E = 29891; F = 45078 },
    @{ Row = 27; E = 29921; F = 45108 }
)

foreach ($item in $firstAndLastDateRows) {
    $ws.Cells.Item($item.Row, 5).Value = $item.E   # Column E = M2_1stDate
    $ws.Cells.Item($item.Row, 6).Value = $item.F   # Column F = M2_LastDate
}

$wb.Save()
